{"js": "// 1) Remove the stray \"yolo\" text from the first paragraph (leave it empty).\nconst yoloResults = context.document.body.search(\"yolo\", { matchCase: true });\nyoloResults.load(\"text\");\nawait context.sync();\nif (yoloResults.items.length > 0) {\n  yoloResults.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Helper: find a unique substring in the whole document body and replace it.\nasync function replaceOnce(searchText, replacement) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"...ces plages horaires, si une plages horaires modifi\u00e9 ou supprim\u00e9 avait \u00e9t\u00e9 r\u00e9serv\u00e9,...\"\n//    -> \"...ces plages horaires, si une plage horaire modifi\u00e9e ou supprim\u00e9e avait \u00e9t\u00e9 r\u00e9serv\u00e9e,...\"\nawait replaceOnce(\n  \"ces plages horaires, si une plages horaires modifi\u00e9 ou supprim\u00e9 avait \u00e9t\u00e9 r\u00e9serv\u00e9,\",\n  \"ces plages horaires, si une plage horaire modifi\u00e9e ou supprim\u00e9e avait \u00e9t\u00e9 r\u00e9serv\u00e9e,\"\n);\n\n// 3) \"(Ou un lieu tier sp\u00e9cifi\u00e9).\" -> \"(Ou un lieu tiers sp\u00e9cifi\u00e9).\"\nawait replaceOnce(\n  \"(Ou un lieu tier sp\u00e9cifi\u00e9).\",\n  \"(Ou un lieu tiers sp\u00e9cifi\u00e9).\"\n);\n\n// 4) \"des pr\u00e9requis,\" -> \"des pr\u00e9-requis,\"\nawait replaceOnce(\n  \"des pr\u00e9requis,\",\n  \"des pr\u00e9-requis,\"\n);\n\n// 5) First \"...pourra \u00eatre configurer.\" (right after \"cours pr\u00e9vus, un syst\u00e8me de rappel par Email\")\nawait replaceOnce(\n  \"cours pr\u00e9vus, un syst\u00e8me de rappel par Email pourra \u00eatre configurer.\",\n  \"cours pr\u00e9vus, un syst\u00e8me de rappel par Email pourra \u00eatre configur\u00e9.\"\n);\n\n// 6) \"diff\u00e9rentes session de cours).\" -> \"diff\u00e9rentes sessions de cours).\"\nawait replaceOnce(\n  \"diff\u00e9rentes session de cours).\",\n  \"diff\u00e9rentes sessions de cours).\"\n);\n\n// 7) Second \"... pourra \u00eatre configurer.\" (after \"d\u2019afficher leur agenda personnel\")\nawait replaceOnce(\n  \"d\u2019afficher leur agenda personnel avec leurs divers cours pr\u00e9vus, un syst\u00e8me de rappel par Email pourra \u00eatre configurer.\",\n  \"d\u2019afficher leur agenda personnel avec leurs divers cours pr\u00e9vus, un syst\u00e8me de rappel par Email pourra \u00eatre configur\u00e9.\"\n);\n\n// 8) Missing space between \"...faire appels aux professeurs\" and \"pour qu\u2019ils donnent...\"\nawait replaceOnce(\n  \"faire appels aux professeurspour qu\u2019ils donnent\",\n  \"faire appels aux professeurs pour qu\u2019ils donnent\"\n);\n", "ps1": "$d = $word.ActiveDocument\n\n# Constants (avoid relying on predefined wd* globals):\n#   wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Text not found: $old\"\n    }\n}\n\n# 1) Remove the stray \"yolo\" text from the first paragraph (leave it empty).\nReplace-Text \"yolo\" \"\"\n\n# 2) \"...ces plages horaires, si une plages horaires modifi\u00e9 ou supprim\u00e9 avait \u00e9t\u00e9 r\u00e9serv\u00e9,...\"\n#    -> \"...ces plages horaires, si une plage horaire modifi\u00e9e ou supprim\u00e9e avait \u00e9t\u00e9 r\u00e9serv\u00e9e,...\"\nReplace-Text \"ces plages horaires, si une plages horaires modifi\u00e9 ou supprim\u00e9 avait \u00e9t\u00e9 r\u00e9serv\u00e9,\" \"ces plages horaires, si une plage horaire modifi\u00e9e ou supprim\u00e9e avait \u00e9t\u00e9 r\u00e9serv\u00e9e,\"\n\n# 3) \"(Ou un lieu tier sp\u00e9cifi\u00e9).\" -> \"(Ou un lieu tiers sp\u00e9cifi\u00e9).\"\nReplace-Text \"(Ou un lieu tier sp\u00e9cifi\u00e9).\" \"(Ou un lieu tiers sp\u00e9cifi\u00e9).\"\n\n# 4) \"des pr\u00e9requis,\" -> \"des pr\u00e9-requis,\"\nReplace-Text \"des pr\u00e9requis,\" \"des pr\u00e9-requis,\"\n\n# 5) & 7) Both \"...pourra \u00eatre configurer.\" occurrences -> \"...pourra \u00eatre configur\u00e9.\"\nReplace-Text \"configurer.\" \"configur\u00e9.\"\n\n# 6) \"diff\u00e9rentes session de cours).\" -> \"diff\u00e9rentes sessions de cours).\"\nReplace-Text \"diff\u00e9rentes session de cours).\" \"diff\u00e9rentes sessions de cours).\"\n\n# 8) Missing space between \"...faire appels aux professeurs\" and \"pour qu'ils donnent...\"\nReplace-Text \"faire appels aux professeurspour qu\u2019ils donnent\" \"faire appels aux professeurs pour qu\u2019ils donnent\"\n"}
